# Replace the 25 three-digit-times-one-digit multiplication problems
# with their new values, per the target diff. Each "old" string occurs
# exactly once in the document, so a single (non "replace-all") find &
# replace per pair is safe. The pair for cell "621x2=1242" -> "429x7=3003"
# is ordered *after* the pair that consumes the original "429x7=3003"
# text (-> "935x2=1870") so that the newly-written "429x7=3003" is not
# immediately re-matched and replaced again.
$d = $word.ActiveDocument

$d.Content.Find.Execute("647×7=4529", $true, $false, $false, $false, $false, $true, 1, $false, "919×9=8271", 2) | Out-Null
$d.Content.Find.Execute("443×6=2658", $true, $false, $false, $false, $false, $true, 1, $false, "999×9=8991", 2) | Out-Null
$d.Content.Find.Execute("116×2=232", $true, $false, $false, $false, $false, $true, 1, $false, "937×6=5622", 2) | Out-Null
$d.Content.Find.Execute("137×7=959", $true, $false, $false, $false, $false, $true, 1, $false, "686×2=1372", 2) | Out-Null
$d.Content.Find.Execute("799×4=3196", $true, $false, $false, $false, $false, $true, 1, $false, "727×9=6543", 2) | Out-Null
$d.Content.Find.Execute("776×6=4656", $true, $false, $false, $false, $false, $true, 1, $false, "371×8=2968", 2) | Out-Null
$d.Content.Find.Execute("191×7=1337", $true, $false, $false, $false, $false, $true, 1, $false, "597×2=1194", 2) | Out-Null
$d.Content.Find.Execute("429×7=3003", $true, $false, $false, $false, $false, $true, 1, $false, "935×2=1870", 2) | Out-Null
$d.Content.Find.Execute("352×8=2816", $true, $false, $false, $false, $false, $true, 1, $false, "686×6=4116", 2) | Out-Null
$d.Content.Find.Execute("731×6=4386", $true, $false, $false, $false, $false, $true, 1, $false, "325×2=650", 2) | Out-Null
$d.Content.Find.Execute("336×9=3024", $true, $false, $false, $false, $false, $true, 1, $false, "417×7=2919", 2) | Out-Null
$d.Content.Find.Execute("177×9=1593", $true, $false, $false, $false, $false, $true, 1, $false, "254×5=1270", 2) | Out-Null
$d.Content.Find.Execute("257×9=2313", $true, $false, $false, $false, $false, $true, 1, $false, "803×9=7227", 2) | Out-Null
$d.Content.Find.Execute("535×6=3210", $true, $false, $false, $false, $false, $true, 1, $false, "973×2=1946", 2) | Out-Null
$d.Content.Find.Execute("674×4=2696", $true, $false, $false, $false, $false, $true, 1, $false, "441×6=2646", 2) | Out-Null
$d.Content.Find.Execute("976×5=4880", $true, $false, $false, $false, $false, $true, 1, $false, "654×6=3924", 2) | Out-Null
$d.Content.Find.Execute("303×2=606", $true, $false, $false, $false, $false, $true, 1, $false, "102×5=510", 2) | Out-Null
$d.Content.Find.Execute("213×8=1704", $true, $false, $false, $false, $false, $true, 1, $false, "136×2=272", 2) | Out-Null
$d.Content.Find.Execute("871×6=5226", $true, $false, $false, $false, $false, $true, 1, $false, "245×2=490", 2) | Out-Null
$d.Content.Find.Execute("544×4=2176", $true, $false, $false, $false, $false, $true, 1, $false, "272×3=816", 2) | Out-Null
$d.Content.Find.Execute("855×4=3420", $true, $false, $false, $false, $false, $true, 1, $false, "753×2=1506", 2) | Out-Null
$d.Content.Find.Execute("303×9=2727", $true, $false, $false, $false, $false, $true, 1, $false, "373×5=1865", 2) | Out-Null
$d.Content.Find.Execute("857×9=7713", $true, $false, $false, $false, $false, $true, 1, $false, "649×7=4543", 2) | Out-Null
$d.Content.Find.Execute("175×5=875", $true, $false, $false, $false, $false, $true, 1, $false, "195×2=390", 2) | Out-Null
$d.Content.Find.Execute("621×2=1242", $true, $false, $false, $false, $false, $true, 1, $false, "429×7=3003", 2) | Out-Null
